$wb = $excel.ActiveWorkbook

# --- Sheet "Entrenadores": insert a new column at K ---------------------
# (shifts existing columns K..U one slot right, to L..V)
$ws = $wb.Worksheets.Item("Entrenadores")
$ws.Columns("K").Insert()

# New header / value for the inserted column (new shared strings)
$ws.Range("K1").Value = "Nombre Foto Carrera Como Jugador"
$ws.Range("K4").Value = "guedecarrerajugador"

# Widen the new column (and its neighbour J, "Nombre Foto Jugador") to fit
# the new "circulo y tamaño a foto jugador" content
$ws.Range("J1").ColumnWidth = 28.33
$ws.Range("K1").ColumnWidth = 28.33

# --- Active sheet / selection -------------------------------------------
# Commit moved the active tab from "Jugadores" to "Entrenadores" and the
# selection there from AB4 to K4 (the new column just added).
$ws.Activate()
$ws.Range("K4").Select()
